$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dict")

# Add the new row of data (row 98) mirroring the style/format of existing rows
$ws.Cells.Item(98, 1).Value = "当月晋升移交客户数(管户)"
$ws.Cells.Item(98, 2).Value = "晋升移交数"

# Match formatting/style of the other cells in column A/B (style index 1 = vertical-center alignment)
$ws.Cells.Item(98, 1).VerticalAlignment = -4108
$ws.Cells.Item(98, 2).VerticalAlignment = -4108

# Update the view: scroll so row 53 is the top-left visible row, and select B98 (A98:B98)
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("B98").Select()
$ws.Range("A98:B98").Select()
